$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths
$ws.Columns.Item(1).ColumnWidth = 18.42578125
$ws.Columns.Item(2).ColumnWidth = 20.42578125
$ws.Columns.Item(3).ColumnWidth = 26.7109375
$ws.Columns.Item(4).ColumnWidth = 20.140625
$ws.Columns.Item(5).ColumnWidth = 20.140625
$ws.Columns.Item(6).ColumnWidth = 29.85546875
$ws.Columns.Item(7).ColumnWidth = 31.42578125

# Row 1 - headers
$ws.Range("A1").Value = "Language"
$ws.Range("B1").Value = "Optimized?"
$ws.Range("C1").Value = "Matrix"
$ws.Range("D1").Value = "Elapsed Nanoseconds"
$ws.Range("E1").Value = "Elapsed Milliseconds"
$ws.Range("F1").Value = "Elapsed Seconds"
$ws.Range("G1").Value = "Transformation performed"
$ws.Range("C1:D1").HorizontalAlignment = -4108

# Row 2
$ws.Range("A2").ClearFormats()
$ws.Range("A2").Value = "C++"
$ws.Range("B2").Value = "No"
$ws.Range("C2").Value = "HarrisInputMatrix.txt"
$ws.Range("C2").HorizontalAlignment = -4152
$ws.Range("D2").Value = 109174000
$ws.Range("D2").NumberFormat = "0.00E+00"
$ws.Range("E2").Formula = "=D2/10^6"
$ws.Range("E2").NumberFormat = "0.00"
$ws.Range("F2").Formula = "=D2/10^9"
$ws.Range("F2").NumberFormat = "General"

# Row 3
$ws.Range("A3").Value = "C++"
$ws.Range("B3").Value = "Yes"
$ws.Range("C3").Value = "HarrisInputMatrix.txt"
$ws.Range("C3").HorizontalAlignment = -4152
$ws.Range("D3").Value = 17294400
$ws.Range("D3").NumberFormat = "0.00E+00"
$ws.Range("E3").Formula = "=D3/10^6"
$ws.Range("E3").NumberFormat = "0.00"
$ws.Range("F3").Formula = "=D3/10^9"
$ws.Range("F3").NumberFormat = "General"
$ws.Range("G3").Value = "OrderStatisticsFiltering: find max"

# Row 4
$ws.Range("A4").Value = "Matlab"
$ws.Range("C4").Value = "HarrisInputMatrix.txt"
$ws.Range("C4").HorizontalAlignment = -4152
$ws.Range("F4").Value = 0.0610010623931884
$ws.Range("D4").Formula = "=F4*10^9"
$ws.Range("D4").NumberFormat = "0.00E+00"
$ws.Range("E4").Formula = "=D4/10^6"
$ws.Range("E4").NumberFormat = "0.00"

# Row 5 - blank but F5 carries a number format
$ws.Range("F5").NumberFormat = "General"

$ws.Range("E7").Select()
